$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "r"
$ws.Range("C4").Value = "g"
$ws.Range("H4").Value = "g"
$ws.Range("F5").Value = "r"
$ws.Range("G5").Value = "r"
$ws.Range("B6").Value = "g"
$ws.Range("E7").Value = "g"
$ws.Range("C9").Value = "r"
$ws.Range("D9").Value = "g"
$ws.Range("F9").Value = "g"
$ws.Range("G9").Value = "g"
$ws.Range("F10").Value = "r"
$ws.Range("G10").Value = "g"
$ws.Range("H10").Value = "g"
